# Add a new "ctvalue" column header in C1, correct a rounded data point in
# C9 (39.8958333333333 -> 40), and leave the selection on the new header
# cell (C1) instead of the previous scroll position / selection (F38).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C1").Value = "ctvalue"
$ws.Range("C9").Value = 40

$ws.Range("C1").Select()
